# Apply the updated inscription/registration counts to the "Inscricoes" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 17
$ws.Range("E17").Value = 117

# Row 36
$ws.Range("E36").Value = 98

# Row 37
$ws.Range("E37").Value = 53
$ws.Range("F37").Value = 30
$ws.Range("H37").Value = 42

# Row 52
$ws.Range("E52").Value = 5

# Row 53
$ws.Range("E53").Value = 6
$ws.Range("F53").Value = 3
$ws.Range("H53").Value = 5

# Row 57
$ws.Range("E57").Value = 15

# Row 63
$ws.Range("E63").Value = 35

# Row 65
$ws.Range("E65").Value = 29

# Row 71
$ws.Range("E71").Value = 35

# Row 89
$ws.Range("E89").Value = 38
